# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted into the "Melón" (Tuna, Segunda)
# series for Vega Monumental Concepción, pushing the existing rows 321-394
# down to 322-395 and extending the sheet's used range to row 395.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 321 - this shifts rows 321:394 down to
# 322:395 (and grows the sheet dimension from R394 to R395) exactly like
# Excel's native "Insert Sheet Rows" command.
$ws.Rows("321").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A321").Value = 11
$ws.Range("B321").Value = "Vega Monumental Concepción"
$ws.Range("C321").Value = "Bíobío"
$ws.Range("D321").Value = 45275
$ws.Range("E321").Value = 8
$ws.Range("F321").Value = 100112027
$ws.Range("G321").Value = "Melón"
$ws.Range("H321").Value = "Tuna"
$ws.Range("I321").Value = "Segunda"
$ws.Range("J321").Value = 1000
$ws.Range("K321").Value = 1500
$ws.Range("L321").Value = 1500
$ws.Range("M321").Value = 1500
$ws.Range("N321").Value = "$/unidad"
$ws.Range("O321").Value = "Región de O'Higgins"
$ws.Range("P321").Value = 1500
$ws.Range("Q321").Value = 1
$ws.Range("R321").Value = "Hortaliza"
